$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.414.34"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.44"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.30"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4488"
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3755"
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07491"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8882"
$ws.Range("E10").Value = "  +4.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.07"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.819.89"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.755"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.89"
$ws.Range("E14").Value = "  +4.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.417"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07113"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008791"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.17"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.423.02"
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.333"
$ws.Range("E22").Value = "  +3.52%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.056.17"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.963"
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.360"
$ws.Range("E26").Value = "  +6.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.52"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.378"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.13"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08872"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7887"
$ws.Range("E32").Value = "  +6.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.201"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.553"
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.925"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.113"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01997"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05334"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.380"
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5327"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1724"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.303"
$ws.Range("E44").Value = "  +19.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.724"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5096"
$ws.Range("E46").Value = "  +5.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.64"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.702"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.79"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06384"
$ws.Range("E51").Value = "  +0.76%  "
